# Add two new "JSON path" columns (N and O) to the Excel-import template,
# matching the header/example rows already used by columns A:M, and move
# the active selection to the new O1 header cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 (example-value row) for the two new columns needs the same
# "wrap text" cell style already used by the other example cells
# (D2:M2) before we populate it, so the new cells share their xf.
$ws.Range("N2:O2").WrapText = $true

# Populate the new cells in the same order the workbook author referenced
# them when they typed the values into the template (bottom example row
# first, then the header row), so new shared-string entries land in the
# same order as the authored file.
$ws.Range("O2").Value = "QnAYesNoBot"
$ws.Range("N2").Value = "Test"
$ws.Range("N1").Value = "clientFilterValues"
$ws.Range("O1").Value = "elicitResponse.responsebot_hook"

# The header row grew a third visual line of text, so it needed more
# height once the new headers were in place.
$ws.Rows(1).RowHeight = 51

# Scroll position / selection: the sheet now opens showing row 1 (no
# frozen/forced topLeftCell override) with the new O1 header selected.
$ws.Range("O1").Select() | Out-Null
